{"js": "// This document has a single table (20 rows x 5 cols); only rows 0, 4, 9, 14, 19\n// (0-based) contain text, 5 cells each. Together with the title paragraph above the\n// table, that is 26 text runs total -- every one of them is replaced below.\nconst body = context.document.body;\n\n// 1) Update the title paragraph (first paragraph of the body).\nconst paragraphs = body.paragraphs;\nparagraphs.load('items');\nawait context.sync();\n\nconst titleRange = paragraphs.items[0].getRange();\ntitleRange.load('text');\nawait context.sync();\n\nif (titleRange.text === '2024-10-05 Saturday') {\n  // insertText(..., replace) on the paragraph's own range swaps only the text\n  // run content, leaving the run's rPr (font/size) and the paragraph's pPr intact.\n  titleRange.insertText('2024-10-06 Sunday', Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Update the table cells. Each entry is [rowIndex, colIndex, oldText, newText]\n//    (0-based row/col, matching Table.getCell). Note '362\u00d76=2172' appears twice\n//    (row 4 col 3, and row 14 col 1) with two different replacements, so cells are\n//    addressed by position rather than by searching for the old text.\nconst tables = body.tables;\ntables.load('items');\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellUpdates = [\n  [0, 0, \"424\u00d72=848\", \"284\u00d77=1988\"],\n  [0, 1, \"811\u00d76=4866\", \"551\u00d74=2204\"],\n  [0, 2, \"748\u00d79=6732\", \"378\u00d74=1512\"],\n  [0, 3, \"128\u00d72=256\", \"168\u00d72=336\"],\n  [0, 4, \"127\u00d76=762\", \"313\u00d77=2191\"],\n  [4, 0, \"235\u00d75=1175\", \"976\u00d73=2928\"],\n  [4, 1, \"563\u00d74=2252\", \"895\u00d79=8055\"],\n  [4, 2, \"491\u00d75=2455\", \"877\u00d79=7893\"],\n  [4, 3, \"362\u00d76=2172\", \"582\u00d77=4074\"],\n  [4, 4, \"984\u00d79=8856\", \"662\u00d78=5296\"],\n  [9, 0, \"563\u00d76=3378\", \"412\u00d77=2884\"],\n  [9, 1, \"313\u00d77=2191\", \"716\u00d77=5012\"],\n  [9, 2, \"228\u00d79=2052\", \"239\u00d77=1673\"],\n  [9, 3, \"769\u00d78=6152\", \"842\u00d78=6736\"],\n  [9, 4, \"289\u00d79=2601\", \"313\u00d75=1565\"],\n  [14, 0, \"904\u00d73=2712\", \"589\u00d73=1767\"],\n  [14, 1, \"362\u00d76=2172\", \"125\u00d79=1125\"],\n  [14, 2, \"233\u00d79=2097\", \"442\u00d79=3978\"],\n  [14, 3, \"246\u00d73=738\", \"466\u00d75=2330\"],\n  [14, 4, \"837\u00d73=2511\", \"869\u00d74=3476\"],\n  [19, 0, \"806\u00d79=7254\", \"648\u00d74=2592\"],\n  [19, 1, \"178\u00d74=712\", \"616\u00d73=1848\"],\n  [19, 2, \"555\u00d73=1665\", \"134\u00d76=804\"],\n  [19, 3, \"919\u00d72=1838\", \"983\u00d76=5898\"],\n  [19, 4, \"559\u00d76=3354\", \"360\u00d79=3240\"]\n];\n\nfor (const [row, col, oldText, newText] of cellUpdates) {\n  const cell = table.getCell(row, col);\n  const cellParagraphs = cell.body.paragraphs;\n  cellParagraphs.load('items');\n  await context.sync();\n\n  const cellRange = cellParagraphs.items[0].getRange();\n  cellRange.load('text');\n  await context.sync();\n\n  if (cellRange.text === oldText) {\n    cellRange.insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# This document has a single table (20 rows x 5 cols); only 1-indexed rows\n# 1, 5, 10, 15, 20 contain text, 5 cells each. Together with the title paragraph\n# above the table, that is 26 text runs total -- every one is replaced below.\n$d = $word.ActiveDocument\n\nfunction Set-RangeText($rng, $oldText, $newText) {\n  # A Range's last character(s) are its paragraph mark / end-of-cell mark; drop\n  # them before comparing/assigning so only the visible text is touched and the\n  # run's formatting (rPr: font/size) plus the paragraph's pPr (jc) survive.\n  $visible = $d.Range($rng.Start, $rng.End - 1)\n  if ($visible.Text -eq $oldText) {\n    $visible.Text = $newText\n  }\n}\n\n# 1) Update the title paragraph (first paragraph of the body).\n$titlePara = $d.Paragraphs.Item(1)\nSet-RangeText $titlePara.Range \"2024-10-05 Saturday\" \"2024-10-06 Sunday\"\n\n# 2) Update the table cells. Each entry is (row, col, oldText, newText) with\n#    1-based row/col, matching Table.Cell(row, col). Note '362\u00d76=2172' appears\n#    twice (row 5 col 4, and row 15 col 2) with two different replacements, so\n#    cells are addressed by position rather than by searching for the old text.\n$table = $d.Tables.Item(1)\n\n$cellUpdates = @(\n  @(1, 1, \"424\u00d72=848\", \"284\u00d77=1988\"),\n  @(1, 2, \"811\u00d76=4866\", \"551\u00d74=2204\"),\n  @(1, 3, \"748\u00d79=6732\", \"378\u00d74=1512\"),\n  @(1, 4, \"128\u00d72=256\", \"168\u00d72=336\"),\n  @(1, 5, \"127\u00d76=762\", \"313\u00d77=2191\"),\n  @(5, 1, \"235\u00d75=1175\", \"976\u00d73=2928\"),\n  @(5, 2, \"563\u00d74=2252\", \"895\u00d79=8055\"),\n  @(5, 3, \"491\u00d75=2455\", \"877\u00d79=7893\"),\n  @(5, 4, \"362\u00d76=2172\", \"582\u00d77=4074\"),\n  @(5, 5, \"984\u00d79=8856\", \"662\u00d78=5296\"),\n  @(10, 1, \"563\u00d76=3378\", \"412\u00d77=2884\"),\n  @(10, 2, \"313\u00d77=2191\", \"716\u00d77=5012\"),\n  @(10, 3, \"228\u00d79=2052\", \"239\u00d77=1673\"),\n  @(10, 4, \"769\u00d78=6152\", \"842\u00d78=6736\"),\n  @(10, 5, \"289\u00d79=2601\", \"313\u00d75=1565\"),\n  @(15, 1, \"904\u00d73=2712\", \"589\u00d73=1767\"),\n  @(15, 2, \"362\u00d76=2172\", \"125\u00d79=1125\"),\n  @(15, 3, \"233\u00d79=2097\", \"442\u00d79=3978\"),\n  @(15, 4, \"246\u00d73=738\", \"466\u00d75=2330\"),\n  @(15, 5, \"837\u00d73=2511\", \"869\u00d74=3476\"),\n  @(20, 1, \"806\u00d79=7254\", \"648\u00d74=2592\"),\n  @(20, 2, \"178\u00d74=712\", \"616\u00d73=1848\"),\n  @(20, 3, \"555\u00d73=1665\", \"134\u00d76=804\"),\n  @(20, 4, \"919\u00d72=1838\", \"983\u00d76=5898\"),\n  @(20, 5, \"559\u00d76=3354\", \"360\u00d79=3240\")\n)\n\nforeach ($u in $cellUpdates) {\n  $row = $u[0]\n  $col = $u[1]\n  $oldText = $u[2]\n  $newText = $u[3]\n  $cell = $table.Cell($row, $col)\n  Set-RangeText $cell.Range $oldText $newText\n}\n"}
